$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = -38750830.2
$ws.Range("P2").Value = -84.0562630957
$ws.Range("Q2").Value = 86646198.45999999
$ws.Range("R2").Value = 187.9483772711
$ws.Range("S2").Value = 28384797.3
$ws.Range("T2").Value = 61.5708096434
$ws.Range("U2").Value = -6631492.88
$ws.Range("V2").Value = -14.3846856277
$ws.Range("Y2").Value = 6631492.88
$ws.Range("Z2").Value = 14.3846856277
$ws.Range("AA2").Value = 92709114.44
$ws.Range("AB2").Value = 201.0997357868
$ws.Range("AC2").Value = 46101062.28
$ws.Range("AD2").ClearContents()
